# Update gh-pages output data (attendance / ticket numbers refreshed at 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览"
$ws1.Range("F4").Value = 3097
$ws1.Range("F5").Value = 1686
$ws1.Range("F6").Value = 2109
$ws1.Range("G7").Value = 218
$ws1.Range("F9").Value = 899
$ws1.Range("F10").Value = 983
$ws1.Range("F11").Value = 223
$ws1.Range("F12").Value = 445
$ws1.Range("F17").Value = 7550
$ws1.Range("F18").Value = 322
$ws1.Range("F20").Value = 204
$ws1.Range("F21").Value = 219
$ws1.Range("F22").Value = 168
$ws1.Range("F24").Value = 522
$ws1.Range("F25").Value = 73
$ws1.Range("F27").Value = 971
$ws1.Range("F29").Value = 1634
$ws1.Range("F31").Value = 1148
$ws1.Range("F34").Value = 29
$ws1.Range("F37").Value = 36
$ws1.Range("F38").Value = 158
$ws1.Range("F39").Value = 317
$ws1.Range("F41").Value = 209

# Sheet "演出"
$ws2.Range("F2").Value = 22

# Sheet "全部类型"
$ws4.Range("F3").Value = 22
$ws4.Range("F7").Value = 3097
$ws4.Range("F8").Value = 1686
$ws4.Range("F9").Value = 2109
$ws4.Range("G10").Value = 218
$ws4.Range("F12").Value = 899
$ws4.Range("F14").Value = 983
$ws4.Range("F15").Value = 223
$ws4.Range("F16").Value = 445
$ws4.Range("F21").Value = 7550
$ws4.Range("F22").Value = 322
$ws4.Range("F25").Value = 204
$ws4.Range("F26").Value = 219
$ws4.Range("F27").Value = 168
$ws4.Range("F29").Value = 522
$ws4.Range("F30").Value = 73
$ws4.Range("F32").Value = 971
$ws4.Range("F34").Value = 1634
$ws4.Range("F36").Value = 1148
$ws4.Range("F39").Value = 29
$ws4.Range("F42").Value = 36
$ws4.Range("F43").Value = 158
$ws4.Range("F44").Value = 317
$ws4.Range("F49").Value = 209
